$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# row 1, col 1: 99÷7= -> 59÷9=
$t.Cell(1, 1).Range.Text = "59÷9="

# row 1, col 2: 48÷3= -> 39÷4=
$t.Cell(1, 2).Range.Text = "39÷4="

# row 1, col 3: 15÷2= -> 89÷7=
$t.Cell(1, 3).Range.Text = "89÷7="

# row 1, col 4: 94÷4= -> 19÷5=
$t.Cell(1, 4).Range.Text = "19÷5="

# row 1, col 5: 47÷3= -> 21÷6=
$t.Cell(1, 5).Range.Text = "21÷6="

# row 5, col 1: 24÷3= -> 63÷4=
$t.Cell(5, 1).Range.Text = "63÷4="

# row 5, col 2: 22÷4= -> 78÷5=
$t.Cell(5, 2).Range.Text = "78÷5="

# row 5, col 3: 47÷9= -> 61÷9=
$t.Cell(5, 3).Range.Text = "61÷9="

# row 5, col 4: 16÷4= -> 57÷5=
$t.Cell(5, 4).Range.Text = "57÷5="

# row 5, col 5: 72÷3= -> 16÷2=
$t.Cell(5, 5).Range.Text = "16÷2="

# row 9, col 1: 79÷9= -> 76÷2=
$t.Cell(9, 1).Range.Text = "76÷2="

# row 9, col 2: 11÷4= -> 34÷2=
$t.Cell(9, 2).Range.Text = "34÷2="

# row 9, col 3: 92÷6= -> 33÷2=
$t.Cell(9, 3).Range.Text = "33÷2="

# row 9, col 4: 86÷8= -> 93÷6=
$t.Cell(9, 4).Range.Text = "93÷6="

# row 9, col 5: 35÷7= -> 83÷8=
$t.Cell(9, 5).Range.Text = "83÷8="

# row 13, col 1: 93÷3= -> 28÷5=
$t.Cell(13, 1).Range.Text = "28÷5="

# row 13, col 2: 69÷2= -> 96÷8=
$t.Cell(13, 2).Range.Text = "96÷8="

# row 13, col 3: 75÷3= -> 75÷4=
$t.Cell(13, 3).Range.Text = "75÷4="

# row 13, col 4: 47÷3= -> 34÷7=
$t.Cell(13, 4).Range.Text = "34÷7="

# row 13, col 5: 33÷4= -> 90÷2=
$t.Cell(13, 5).Range.Text = "90÷2="

# row 17, col 1: 39÷3= -> 83÷3=
$t.Cell(17, 1).Range.Text = "83÷3="

# row 17, col 2: 34÷6= -> 87÷9=
$t.Cell(17, 2).Range.Text = "87÷9="

# row 17, col 3: 88÷7= -> 22÷2=
$t.Cell(17, 3).Range.Text = "22÷2="

# row 17, col 4: 51÷3= -> 35÷5=
$t.Cell(17, 4).Range.Text = "35÷5="

# row 17, col 5: 95÷6= -> 14÷3=
$t.Cell(17, 5).Range.Text = "14÷3="
